$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text helper pattern: prefix value with a leading apostrophe so the
# COM layer stores it as text (matching the source inlineStr cells) instead
# of auto-converting number-like strings (e.g. "1.00") to numeric values.
# Resetting the Style back to "Normal" afterwards clears the quote-prefix
# cell style Excel applies, keeping cell formatting identical to the original.

$ws.Range("D2").Value = "'60.584.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.28%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.660.03"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'568.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +6.00%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'145.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.48%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.11%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +6.92%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'6.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.84%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +4.16%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E12").Value = "'  +2.25%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'3.120.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.81%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'60.519.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.27%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'21.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.09%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'2.649.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.19%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +2.94%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +3.55%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'343.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.65%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +2.22%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'6.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.35%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.73%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.14%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'66.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.83%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +5.87%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +1.35%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.30%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +1.97%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +4.14%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.06%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.27%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'6.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.44%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'156.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.64%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +2.15%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +4.39%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "'Fetch.AI"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.910"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +11.62%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'SuiNetwork"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.903"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.46%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +5.15%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'37.50"
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = "'  +5.41%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'304.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +7.64%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +2.03%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +0.09%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +0.50%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0977"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.40%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +2.87%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'19.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.13%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'10.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.51%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'Aave"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'125.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +11.81%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'VeChain"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0235"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.38%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'Maker"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'1.968.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.93%  "
$ws.Range("E51").Style = "Normal"
